# Generate Report for Handoff
# Adds two new files (6f5ba720-... and b577a97b-...) into the localization
# status report, and refreshes the "Latest Handoff" status/date for every
# file row on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Row data, in final top-to-bottom order, for every sheet.
# ---------------------------------------------------------------------
$guids = @(
    "6f5ba720-2541-488a-af3e-5fa682e96881",
    "81da33ca-0519-4365-b4c1-dc94f41f1156",
    "b577a97b-1ff7-420c-aaba-249f9ebd3035",
    "f983ba90-87a7-437f-b6c0-1367667f4d87"
)
$hashes = @(
    "eb182456c8705f50c5bcba5495af8f3cf753368d",
    "5c47de6e1c32b19976a262704a561d1b1bd2ec41",
    "71571e19be4f5a65b9a2fae21bf2ad59557bbf1f",
    "9b0be5b62d3e0de653f5e4eb9e2c52deac87cdde"
)

$mdBase  = "https://github.com/OpenLocalizationTest/oltest/blob/f8f40194a6bd795944a365ab7b6155b40eb3bc92/e2e/"
$zhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d06916dd49e381526c17576fe5a3f4fcd4c1242/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a2b0d25a2e97f54b7a7a265bb0a363d0c6be0fe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$zhHandoffDate = "2016-03-11 06:26:44"
$deHandoffDate = "2016-03-11 06:26:52"
$overviewDate  = "2016-26-11 06:26:52"
$epoch         = "0001-01-01 00:00:00"
$status        = "Ready for handoff"

# ===========================================================================
# Sheet "Overview": A=File name, B=zh-cn, C=de-de, D=Latest Handoff Date
# ===========================================================================
$wsOv = $wb.Worksheets.Item("Overview")

# Make room: rows 2 & 3 already exist (81da33ca, f983ba90) - insert one
# blank row above each of them so the new guids can be slotted in, giving
# final order: 6f5ba720, 81da33ca, b577a97b, f983ba90
$wsOv.Range("A2").EntireRow.Insert()
$wsOv.Range("A4").EntireRow.Insert()

for ($i = 0; $i -lt $guids.Length; $i++) {
    $r = $i + 2
    $g = $guids[$i]
    $wsOv.Cells.Item($r, 2).Value = $status
    $wsOv.Cells.Item($r, 3).Value = $status
    $wsOv.Cells.Item($r, 4).Value = $overviewDate
}

$wsOv.Hyperlinks.Delete()
for ($i = 0; $i -lt $guids.Length; $i++) {
    $r = $i + 2
    $g = $guids[$i]
    $wsOv.Hyperlinks.Add($wsOv.Cells.Item($r, 1), ($mdBase + $g + ".md"), "", "", ($g + ".md"))
}

# ===========================================================================
# Sheets "zh-cn" / "de-de": both share the same 11-column layout
#  A=Source File Name  B=File Extension  C=Status  D=Latest Handoff File
#  E=Latest Handoff Datetime  H=Latest Handback DateTime  I=Handoff Reason
# ===========================================================================
function Update-LangSheet($ws, $xlfBase, $handoffDate) {
    # Existing rows 2 (81da33ca) & 3 (f983ba90) - insert blank rows above
    # each so the new guids fit in, same final order as Overview.
    $ws.Range("A2").EntireRow.Insert()
    $ws.Range("A4").EntireRow.Insert()

    $ws.Hyperlinks.Delete()

    for ($i = 0; $i -lt $guids.Length; $i++) {
        $r = $i + 2
        $g = $guids[$i]
        $h = $hashes[$i]
        $xlfName = $g + "." + $h + "." + $xlfBase.Suffix

        $ws.Cells.Item($r, 1).Value = $g + ".md"
        $ws.Cells.Item($r, 2).Value = ".md"
        $ws.Cells.Item($r, 3).Value = $status
        $ws.Cells.Item($r, 4).Value = $xlfName
        $ws.Cells.Item($r, 5).Value = $handoffDate
        $ws.Cells.Item($r, 8).Value = $epoch
        $ws.Cells.Item($r, 9).Value = "Include"

        $ws.Hyperlinks.Add($ws.Cells.Item($r, 1), ($mdBase + $g + ".md"), "", "", ($g + ".md"))
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), ($mdBase + $g + ".md"), "", "", ".md")
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 4), ($xlfBase.Url + $xlfName), "", "", $xlfName)
    }
}

$wsZh = $wb.Worksheets.Item("zh-cn")
$zh = @{ Url = $zhBase; Suffix = "zh-cn.xlf" }
Update-LangSheet $wsZh $zh $zhHandoffDate

$wsDe = $wb.Worksheets.Item("de-de")
$de = @{ Url = $deBase; Suffix = "de-de.xlf" }
Update-LangSheet $wsDe $de $deHandoffDate

Write-Host "Report updated: added" $guids.Length "rows total (2 new files) across Overview/zh-cn/de-de."
